{"js": "// Word JS API (Office.js) edit script.\n// Body is `async (context) => { ... }`.\n\nconst body = context.document.body;\n\n// --- 1. Fix the duplicated \"librari .../labora ...\" typo fragments inside\n//        the intro paragraph: \"librari libraries,labora laboratory\" ->\n//        \"libraries,laboratory\".\nconst typoResults = body.search(\"librari libraries,labora laboratory\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"libraries,laboratory\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Reword \"...Among of the society that will be solved...\" to\n//        \"...Among of  society problem that will be solved...\".\nconst amongResults = body.search(\"the society\", { matchCase: true });\namongResults.load(\"text\");\nawait context.sync();\n\nif (amongResults.items.length > 0) {\n  amongResults.items[0].insertText(\" society problem\", \"Replace\");\n  await context.sync();\n}\n\n// --- 3. Append four new paragraphs right after the\n//        \"4.circumstances to produce sufficient witnesses.\" paragraph\n//        (i.e. right before the trailing empty paragraph).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"4.circumstances to produce sufficient witnesses.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  // Fallback: last non-empty paragraph before the trailing blank one.\n  anchor = paragraphs.items[paragraphs.items.length - 2];\n}\n\nconst newParaTexts = [\n  \"Now , this project will more  increase safety at sensitive areas and will also reduce the number of crimes at the different areas.\",\n  \"  Prepared by\",\n  \"  1.godfrey mwakilembe-20100523140078\",\n  \"  2.nelson warioba-20100523140080\",\n];\n\nlet last = anchor;\nfor (const text of newParaTexts) {\n  last = last.insertParagraph(text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix the duplicated \"librari .../labora ...\" typo fragments inside\n#        the intro paragraph: \"librari libraries,labora laboratory\" ->\n#        \"libraries,laboratory\".\n$d.Content.Find.Execute(\n    \"librari libraries,labora laboratory\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"libraries,laboratory\",\n    2\n) | Out-Null\n\n# --- 2. Reword \"...Among of the society that will be solved...\" to\n#        \"...Among of  society problem that will be solved...\".\n$d.Content.Find.Execute(\n    \"the society\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \" society problem\",\n    2\n) | Out-Null\n\n# --- 3. Append four new paragraphs right after the\n#        \"4.circumstances to produce sufficient witnesses.\" paragraph\n#        (i.e. right before the trailing empty paragraph).\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"*4.circumstances to produce sufficient witnesses.*\") {\n        $anchor = $candidate\n        break\n    }\n}\nif ($anchor -eq $null) {\n    # Fallback: second-to-last paragraph (just before the trailing blank one).\n    $anchor = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n}\n\n$newParaTexts = @(\n    \"Now , this project will more  increase safety at sensitive areas and will also reduce the number of crimes at the different areas.\",\n    \"  Prepared by\",\n    \"  1.godfrey mwakilembe-20100523140078\",\n    \"  2.nelson warioba-20100523140080\"\n)\n\nforeach ($text in $newParaTexts) {\n    $anchor.Range.InsertParagraphAfter()\n    $anchor = $anchor.Next()\n    $anchor.Range.Text = $text\n}\n"}
